# "Added Course Materials - Day 4"
#
# Removes a handful of leftover decorative rectangle shapes (plain
# custGeom boxes with no text) from slides 7, 8 and 11, and nudges the
# page-number "Text Box 27" on slide 7 into its new position.

$p = $ppt.ActivePresentation

# --- Slide 7 --------------------------------------------------------
$s7 = $p.Slides.Item(7)

# Remove the eight stray decorative boxes ("object 15".."object 22").
foreach ($name in @("object 15","object 16","object 17","object 18","object 19","object 20","object 21","object 22")) {
    $s7.Shapes.Item($name).Delete()
}

# Reposition the page-number text box.
$tb = $s7.Shapes.Item("Text Box 27")
$tb.Left = 676.35001
$tb.Top = 420.0

# --- Slide 8 ----------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item("object 15").Delete()

# --- Slide 11 -----------------------------------------------------------
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item("object 21").Delete()
